# Updates cryptos list data (prices and 1h volume change) to reflect latest scrape.
# A leading apostrophe forces Excel to store the value as text (matching the
# original inlineStr cell type) instead of auto-converting numeric-looking strings
# such as "244.54" into numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'36.302.65"
$ws.Range("E2").Value = "'  -3.32%  "

$ws.Range("D3").Value = "'1.976.09"
$ws.Range("E3").Value = "'  -4.03%  "

$ws.Range("E4").Value = "'  +0.06%  "

$ws.Range("D5").Value = "'244.54"
$ws.Range("E5").Value = "'  -3.80%  "

$ws.Range("E6").Value = "'  -4.64%  "

$ws.Range("D7").Value = "'58.57"
$ws.Range("E7").Value = "'  -12.74%  "

$ws.Range("E8").Value = "'  +0.01%  "

$ws.Range("E9").Value = "'  -6.28%  "

$ws.Range("D10").Value = "'56.86"
$ws.Range("E10").Value = "'  -5.27%  "

$ws.Range("D11").Value = "'0.0839"
$ws.Range("E11").Value = "'  +8.64%  "

$ws.Range("E12").Value = "'  -0.76%  "

$ws.Range("D13").Value = "'22.97"
$ws.Range("E13").Value = "'  -3.44%  "

$ws.Range("D14").Value = "'0.854"
$ws.Range("E14").Value = "'  -9.07%  "

$ws.Range("D15").Value = "'2.267.30"
$ws.Range("E15").Value = "'  -3.91%  "

$ws.Range("E16").Value = "'  -7.64%  "

$ws.Range("D17").Value = "'5.41"
$ws.Range("E17").Value = "'  -6.13%  "

$ws.Range("D18").Value = "'1.976.13"
$ws.Range("E18").Value = "'  -4.01%  "

$ws.Range("D19").Value = "'36.165.34"
$ws.Range("E19").Value = "'  -3.45%  "

$ws.Range("D20").Value = "'0.0₃0881"
$ws.Range("E20").Value = "'  -0.05%  "

$ws.Range("D21").Value = "'70.19"
$ws.Range("E21").Value = "'  -4.90%  "

$ws.Range("D22").Value = "'5.26"
$ws.Range("E22").Value = "'  -4.64%  "

$ws.Range("D23").Value = "'233.59"
$ws.Range("E23").Value = "'  -3.06%  "

$ws.Range("E24").Value = "'  +0.06%  "

$ws.Range("D25").Value = "'2.51"
$ws.Range("E25").Value = "'  -6.53%  "

$ws.Range("D26").Value = "'2.28"
$ws.Range("E26").Value = "'  -6.96%  "

$ws.Range("D27").Value = "'9.86"
$ws.Range("E27").Value = "'  -2.63%  "

$ws.Range("D28").Value = "'162.97"
$ws.Range("E28").Value = "'  +0.20%  "

$ws.Range("E29").Value = "'  -2.33%  "

$ws.Range("D30").Value = "'0.132"
$ws.Range("E30").Value = "'  -5.03%  "

$ws.Range("D31").Value = "'0.120"
$ws.Range("E31").Value = "'  -2.76%  "

$ws.Range("E32").Value = "'  -4.42%  "

$ws.Range("B33").Value = "'Hedera"
$ws.Range("C33").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.0681"
$ws.Range("E33").Value = "'  +6.85%  "

$ws.Range("B34").Value = "'Filecoin"
$ws.Range("C34").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.85"
$ws.Range("E34").Value = "'  -7.41%  "

$ws.Range("E35").Value = "'  -7.21%  "

$ws.Range("E36").Value = "'  -2.64%  "

$ws.Range("E37").Value = "'  +0.04%  "

$ws.Range("E38").Value = "'  -0.91%  "

$ws.Range("E39").Value = "'  -8.91%  "

$ws.Range("D40").Value = "'2.95"
$ws.Range("E40").Value = "'  -6.66%  "

$ws.Range("E41").Value = "'  -4.15%  "

$ws.Range("D42").Value = "'0.0954"
$ws.Range("E42").Value = "'  -8.19%  "

$ws.Range("E43").Value = "'  -5.63%  "

$ws.Range("D44").Value = "'0.0213"
$ws.Range("E44").Value = "'  -3.80%  "

$ws.Range("E45").Value = "'  -6.45%  "

$ws.Range("D46").Value = "'91.63"
$ws.Range("E46").Value = "'  -5.89%  "

$ws.Range("D47").Value = "'16.05"
$ws.Range("E47").Value = "'  -13.12%  "

$ws.Range("E48").Value = "'  -7.20%  "

$ws.Range("D49").Value = "'1.358.23"
$ws.Range("E49").Value = "'  -4.31%  "

$ws.Range("E50").Value = "'  -4.69%  "

$ws.Range("D51").Value = "'44.75"
$ws.Range("E51").Value = "'  -7.31%  "
